# Add a new "12-10-2020" column (AA) to the deceased cases sheet.
# Column AA mirrors the formatting of column Z:
#   - AA1 is a header cell styled like Z1 (bold, bordered) holding the date as text.
#   - AA2:AA36 are plain numeric cells (same style as Z2:Z36, i.e. default style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell: match the look of Z1 (bold font, thin box border, centered/top
# aligned) and store the date as literal text (like the other date headers)
# instead of letting it be auto-parsed into a date serial number.
$aa1 = $ws.Range("AA1")
$aa1.NumberFormat = "@"
$aa1.Font.Bold = $true
$aa1.Font.Name = "Calibri"
$aa1.Font.Size = 11
$aa1.HorizontalAlignment = -4108  # xlCenter
$aa1.VerticalAlignment = -4160    # xlTop
$aa1.Borders.LineStyle = 1        # xlContinuous (thin box border)
$aa1.Value2 = "12-10-2020"

# Data values for the new column, row by row.
$aaValues = @{
    2  = 55
    3  = 6224
    4  = 24
    5  = 816
    6  = 946
    7  = 191
    8  = 1253
    9  = 2
    10 = 5769
    11 = 507
    12 = 3566
    13 = 1579
    14 = 250
    15 = 1322
    16 = 787
    17 = 9966
    18 = 1003
    19 = 64
    20 = 2624
    21 = 40349
    22 = 91
    23 = 63
    24 = 0
    25 = 17
    26 = 1022
    27 = 563
    28 = 3833
    29 = 1650
    30 = 55
    31 = 10252
    32 = 1228
    33 = 316
    34 = 747
    35 = 6394
    36 = 5622
}

foreach ($row in $aaValues.Keys) {
    $ws.Cells.Item($row, 27).Value = $aaValues[$row]
}
